$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "TheWitchofDelusions"
$ws.Range("C36").Value = "The Witch of Delusions"
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = "Delusional Memories"
$ws.Range("I36").Value = 384
$ws.Range("J36").Value = 80
